# Apply crypto price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number must be forced to
# Text so they stay consistent with the rest of the (text-formatted) column.
$textCells = @('D5', 'D6', 'D8', 'D18', 'D20', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D34', 'D36', 'D37', 'D38', 'D41', 'D42', 'D43', 'D44', 'D45', 'D47', 'D48', 'D51')
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range('D2').Value = '35.137.21'
$ws.Range('E2').Value = '  +1.37%  '
$ws.Range('D3').Value = '1.833.12'
$ws.Range('E3').Value = '  +1.37%  '
$ws.Range('E4').Value = '  +0.53%  '
$ws.Range('D5').Value = '232.72'
$ws.Range('E5').Value = '  +3.36%  '
$ws.Range('D6').Value = '0.617'
$ws.Range('E6').Value = '  +2.14%  '
$ws.Range('E7').Value = '  +0.61%  '
$ws.Range('D8').Value = '42.93'
$ws.Range('E8').Value = '  +5.86%  '
$ws.Range('E9').Value = '  +6.18%  '
$ws.Range('E10').Value = '  +2.66%  '
$ws.Range('E11').Value = '  +0.50%  '
$ws.Range('D12').Value = '2.099.73'
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').Value = '1.822.47'
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('E14').Value = '  +2.77%  '
$ws.Range('E15').Value = '  +4.80%  '
$ws.Range('E16').Value = '  +6.72%  '
$ws.Range('D17').Value = '35.117.06'
$ws.Range('E17').Value = '  +1.23%  '
$ws.Range('D18').Value = '70.35'
$ws.Range('E18').Value = '  +3.68%  '
$ws.Range('E19').Value = '  +2.96%  '
$ws.Range('D20').Value = '240.47'
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('E21').Value = '  +7.21%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '4.59'
$ws.Range('E22').Value = '  +12.14%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '1.01'
$ws.Range('E23').Value = '  +0.55%  '
$ws.Range('D24').Value = '2.24'
$ws.Range('E24').Value = '  +3.95%  '
$ws.Range('D25').Value = '171.52'
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('D26').Value = '7.83'
$ws.Range('E26').Value = '  +1.57%  '
$ws.Range('D27').Value = '17.56'
$ws.Range('E27').Value = '  +0.48%  '
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('E29').Value = '  +29.59%  '
$ws.Range('E30').Value = '  +0.58%  '
$ws.Range('D31').Value = '3.352.07'
$ws.Range('E31').Value = '  +37.96%  '
$ws.Range('E32').Value = '  +7.70%  '
$ws.Range('E33').Value = '  +3.78%  '
$ws.Range('D34').Value = '4.01'
$ws.Range('E34').Value = '  +4.13%  '
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').Value = '93.65'
$ws.Range('E36').Value = '  +10.49%  '
$ws.Range('D37').Value = '0.680'
$ws.Range('E37').Value = '  +5.76%  '
$ws.Range('D38').Value = '1.11'
$ws.Range('E38').Value = '  +5.39%  '
$ws.Range('D39').Value = '1.326.98'
$ws.Range('E39').Value = '  +1.14%  '
$ws.Range('E40').Value = '  +2.97%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D41').Value = '1.27'
$ws.Range('E41').Value = '  +2.11%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '0.997'
$ws.Range('E42').Value = '  +5.55%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '2.37'
$ws.Range('E43').Value = '  +0.65%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '14.94'
$ws.Range('E44').Value = '  -0.51%  '
$ws.Range('D45').Value = '2.46'
$ws.Range('E45').Value = '  +1.16%  '
$ws.Range('E46').Value = '  -0.18%  '
$ws.Range('D47').Value = '6.24'
$ws.Range('E47').Value = '  +8.70%  '
$ws.Range('D48').Value = '0.0508'
$ws.Range('E48').Value = '  -2.53%  '
$ws.Range('D49').Value = '2.008.25'
$ws.Range('E49').Value = '  +2.01%  '
$ws.Range('E50').Value = '  +0.70%  '
$ws.Range('D51').Value = '101.03'
$ws.Range('E51').Value = '  -0.17%  '
